$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.946.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8207'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3170'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.58'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06988'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08020'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7510'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.898.84'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.201'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.947.78'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.884'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007759'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.24%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.154.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9989'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9993'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.967'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1600'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +24.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.229'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.092'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.365'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.512'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.300'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05569'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.53%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.091'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.268'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7338'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.708'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01922'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.785'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4419'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.995'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9987'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8369'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.31%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.891'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.29%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.579'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.715'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '987.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.060.23'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.03%  '
